$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accession/x/y/z records appended to the validation set (rows 23-41)
$data = @(
    @(23, 9056633, 30, 57, 26),
    @(24, 9038021, 25, 31, 23),
    @(25, 9028403, 21, 42, 20),
    @(26, 8977015, 31, 57, 23),
    @(27, 8968387, 32, 54, 19),
    @(28, 8920704, 32, 54, 26),
    @(29, 8920679, 26, 54, 28),
    @(30, 8840552, 24, 46, 21),
    @(31, 8781623, 29, 35, 20),
    @(32, 8765787, 29, 46, 18),
    @(33, 8602870, 38, 48, 21),
    @(34, 8602775, 24, 47, 20),
    @(35, 8709670, 28, 64, 23),
    @(36, 8699113, 32, 46, 24),
    @(37, 8666330, 26, 45, 20),
    @(38, 8665805, 31, 47, 21),
    @(39, 8665616, 29, 38, 19),
    @(40, 8662841, 27, 48, 21),
    @(41, 8636498, 36, 46, 21)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Leave the active selection on the last entered cell, matching the author's edit
$ws.Range("D41").Select()
